# AOCS sizing model complete
# Applies: new "ThrustVectors" sheet, TTC D3 value change, and the
# various selection / active-tab bookkeeping changes recorded in the diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. TTC sheet: offset1 for TTC-nadir (D3) goes from 0.25 -> 0.4
# ---------------------------------------------------------------------
$ttc = $wb.Worksheets.Item("TTC")
$ttc.Range("D3").Value = 0.4

# ---------------------------------------------------------------------
# 2. Add the new "ThrustVectors" sheet after "Thrusters" (last sheet)
# ---------------------------------------------------------------------
$thrusters = $wb.Worksheets.Item("Thrusters")
$tv = $wb.Worksheets.Add($null, $thrusters)
$tv.Name = "ThrustVectors"

# Header row
$tv.Cells.Item(1,1).Value = "name"
$tv.Cells.Item(1,2).Value = "face1"
$tv.Cells.Item(1,3).Value = "face2"
$tv.Cells.Item(1,4).Value = "x"
$tv.Cells.Item(1,5).Value = "y"
$tv.Cells.Item(1,6).Value = "z"
$tv.Cells.Item(1,7).Value = "pair"

# Data rows: name, face1, face2, x, y, z, pair
$rows = @(
    @("att1",  "x+", "y-", -1,  0,  0, "att4"),
    @("att2",  "x+", "y+", -1,  0,  0, "att3"),
    @("att3",  "x-", "y-",  1,  0,  0, "att2"),
    @("att4",  "x-", "y+",  1,  0,  0, "att1"),
    @("att5",  "y+", "z+",  0, -1,  0, "att8"),
    @("att6",  "y+", "z-",  0, -1,  0, "att7"),
    @("att7",  "y-", "z+",  0,  1,  0, "att6"),
    @("att8",  "y-", "z-",  0,  1,  0, "att5"),
    @("att9",  "z+", "x+",  0,  0, -1, "att12"),
    @("att10", "z+", "x-",  0,  0, -1, "att11"),
    @("att11", "z-", "x+",  0,  0,  1, "att10"),
    @("att12", "z-", "x-",  0,  0,  1, "att9")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $tv.Cells.Item($r,1).Value = $row[0]
    $tv.Cells.Item($r,2).Value = $row[1]
    $tv.Cells.Item($r,3).Value = $row[2]
    $tv.Cells.Item($r,4).Value = $row[3]
    $tv.Cells.Item($r,5).Value = $row[4]
    $tv.Cells.Item($r,6).Value = $row[5]
    $tv.Cells.Item($r,7).Value = $row[6]
}

# ---------------------------------------------------------------------
# 3. Selections on each sheet (restores the cursor position recorded
#    for every sheetView) and tab/active-sheet bookkeeping. The sheet
#    selected last becomes the active tab, so "Calculations Rough
#    Input" is selected last to match the target workbookView
#    (no activeTab override, i.e. tab 0 active).
# ---------------------------------------------------------------------
$tv.Range("A6").Select() | Out-Null

$eps = $wb.Worksheets.Item("EPS")
$eps.Range("F2").Select() | Out-Null

$thrusters.Range("H27").Select() | Out-Null

$ttc.Range("D3").Select() | Out-Null

$calc = $wb.Worksheets.Item("Calculations Rough Input")
$calc.Range("K8").Select() | Out-Null
